$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(39583, 2008, $null, 2009, 0.07975464681371225),
    @(39765, 2008, $null, 2009, -4.700509864312973),
    @(39948, 2009, -3.017636378339217, 2010, -3.229247082222797),
    @(40130, 2009, -1.31761954538927, 2010, -0.01655958389530365),
    @(40310, 2010, 0.003352386816724007, 2011, -1.461031976610316),
    @(40494, 2010, 2.771597318554297, 2011, 3.579142225970444),
    @(40676, 2011, 3.925837669383347, 2012, 2.641604203902781),
    @(40862, 2011, 1.799362536952542, 2012, -0.289184878867832),
    @(41044, 2012, 0.2381541440396262, 2013, 1.60268309892857),
    @(41228, 2012, 2.123182427147152, 2013, 5.963492031746176),
    @(41409, 2013, 4.993892964711621, 2014, 2.260118192030736),
    @(41592, 2013, 7.317297369612819, 2014, 7.523777575896196),
    @(41774, 2014, 6.711795724673664, 2015, 6.409878804372982),
    @(41957, 2014, 4.260319658857736, 2015, 2.532215190177589),
    @(42137, 2015, 0.5121603413743347, 2016, 3.290935868252554),
    @(42321, 2015, -0.05262415810141086, 2016, 2.051185924063259),
    @(42503, 2016, 1.745565778643887, 2017, 0.7985845180024986),
    @(42689, 2016, 1.459778471779982, 2017, 0.4575538530338541),
    @(42867, 2017, 2.687500891103922, 2018, 1.922191950024699),
    @(43053, 2017, 3.002208343813528, 2018, 2.600569166164624),
    @(43145, 2018, 3.339205815020496, 2019, 2.646788941483735),
    @(43235, 2018, 3.654655474034474, 2019, 3.068403604789749),
    @(43326, 2018, 4.036117574265741, 2019, 3.628675245064317),
    @(43418, 2018, 4.020433260014977, 2019, 3.605726003451304),
    @(43510, 2019, 3.813466308501412, 2020, 3.884173085820986),
    @(43600, 2019, 3.712036718632117, 2020, 3.908921577463587),
    @(43691, 2019, 3.551357200054261, 2020, 3.524277826276134),
    @(43783, 2019, 3.53224976671227, 2020, 3.490656491795074),
    @(43875, 2020, 2.732790977059629, 2021, 3.412957258051663),
    @(43966, 2020, 2.849400388885992, 2021, 3.535456592693387),
    @(44068, 2020, -0.985458715495402, 2021, -3.099271113627677),
    @(44159, 2020, -0.985458715495402, 2021, -2.347097924577757),
    @(44251, 2021, -5.285337128797329, 2022, -2.69436300797079),
    @(44341, 2021, -4.741003096464214, 2022, -2.156795995006056),
    @(44432, 2021, -4.365687260408224, 2022, -0.9191921099315992),
    @(44525, 2021, -4.365687260408224, 2022, -0.1803381976702711),
    @(44617, 2022, 1.871837441670499, 2023, -0.8573220257725445),
    @(44706, 2022, 1.194925448553708, 2023, -2.068675356622807),
    @(44798, 2022, 1.27347919322387, 2023, -2.137023292796481),
    @(44890, 2022, 1.27347919322387, 2023, -1.152671696465724),
    @(44981, 2023, -0.507642354784088, 2024, 0.07646803116447831),
    @(45071, 2023, -1.084365158506884, 2024, -1.089896342664354),
    @(45163, 2023, -1.339436245206127, 2024, -1.277727682704721),
    @(45254, 2023, -1.339436245206127, 2024, -2.785556326028149),
    @(45345, 2024, -2.87408779878463, 2025, -1.155307395925487),
    @(45436, 2024, -3.40787540386569, 2025, -2.191935020614488),
    @(45534, 2024, -3.451527003230626, 2025, -2.998472503487815),
    @(45618, 2024, -3.451527003230626, 2025, -2.452009576682213),
    @(45713, 2025, -2.164008261629446, 2026, -2.337170009804157),
    @(45800, 2025, -1.853660925652212, 2026, -1.878672029998096),
    @(45891, 2025, -1.75044229618867, 2026, -1.760724207457021)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    if ($row[2] -ne $null) {
        $ws.Cells.Item($r, 3).Value = $row[2]
    } else {
        $ws.Cells.Item($r, 3).ClearContents()
    }
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}

# Remove old row 53 (last row is now 52)
$ws.Rows.Item(53).Delete()

Write-Host "Done. UsedRange: $($ws.UsedRange.Address())"
